# The two species-observation records on rows 17 and 18 were swapped
# (re-ordered). Apply the swap by writing each changed field directly to
# its final value rather than re-ordering whole rows, so cell typing
# (numbers vs. text) is fully under our control.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as TEXT (Excel would otherwise
# auto-convert strings that merely look numeric/date-like, e.g. "2" or
# "2023-09-16", into numbers/dates on assignment). Restores the cell's
# original style afterwards so no stray number-format gets left behind.
function Set-TextValue {
    param($Cell, $Text)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value2 = $Text
    $Cell.Style = $origStyle
}

# Helper: make sure a cell exists but is empty (used where the "other"
# row had a blank cell present in that column). Briefly forcing text mode
# and then clearing is what makes Excel actually materialise the blank
# cell instead of leaving it absent from the sheet.
function Set-EmptyCell {
    param($Cell)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value2 = "x"
    $Cell.Value2 = ""
    $Cell.Style = $origStyle
}

# ----- Row 17 becomes what row 18 used to be -----
$ws.Cells.Item(17, 1).Value2 = 112128590        # A - Id
$ws.Cells.Item(17, 2).Value2 = 56446             # B - Taxonsorteringsordning
$ws.Cells.Item(17, 4).Value2 = "NT"              # D - Rödlistade
$ws.Cells.Item(17, 5).Value2 = 100049            # E - TaxonId
$ws.Cells.Item(17, 6).Value2 = "Spillkråka"      # F - Artnamn
$ws.Cells.Item(17, 7).Value2 = "Dryocopus martius"  # G - Vetenskapligt namn
$ws.Cells.Item(17, 8).Value2 = "(Linnaeus, 1758)"   # H - Auktor
Set-TextValue $ws.Cells.Item(17, 9) "1"          # I - Antal
$ws.Cells.Item(17, 10).ClearContents()           # J - Enhet (now blank)
Set-EmptyCell $ws.Cells.Item(17, 12)             # L - Kön (now present & blank)
$ws.Cells.Item(17, 13).Value2 = "färska spår"    # M - Aktivitet
$ws.Cells.Item(17, 16).Value2 = "Persbo gruva, Upl"  # P - Lokalnamn
$ws.Cells.Item(17, 17).Value2 = 654853           # Q - Ost
$ws.Cells.Item(17, 18).Value2 = 6675824          # R - Nord
Set-TextValue $ws.Cells.Item(17, 26) "12:44"     # Z - Starttid
Set-TextValue $ws.Cells.Item(17, 28) "12:44"     # AB - Sluttid

# ----- Row 18 becomes what row 17 used to be -----
$ws.Cells.Item(18, 1).Value2 = 112129050
$ws.Cells.Item(18, 2).Value2 = 90835
$ws.Cells.Item(18, 4).Value2 = "LC"
$ws.Cells.Item(18, 5).Value2 = 5964
$ws.Cells.Item(18, 6).Value2 = "Fjällig taggsvamp s.str."
$ws.Cells.Item(18, 7).Value2 = "Sarcodon imbricatus s.str."
$ws.Cells.Item(18, 8).Value2 = "(L.:Fr.) P.Karst."
Set-TextValue $ws.Cells.Item(18, 9) "2"
$ws.Cells.Item(18, 10).Value2 = "fruktkroppar"    # J - now present
$ws.Cells.Item(18, 12).ClearContents()            # L - no longer present
$ws.Cells.Item(18, 13).ClearContents()            # M - no longer present
$ws.Cells.Item(18, 16).Value2 = "Persbomossen, Upl"
$ws.Cells.Item(18, 17).Value2 = 654924
$ws.Cells.Item(18, 18).Value2 = 6675762
Set-TextValue $ws.Cells.Item(18, 26) "12:20"
Set-TextValue $ws.Cells.Item(18, 28) "12:20"
